$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4

# Update the selection to D3, matching the diff
$ws.Range("D3").Select()
